# Update the "want to go" count (column F) for a handful of events in the
# "展览" and "全部类型" worksheets. This mirrors the gh-pages data refresh
# captured in the commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (index/row -> new value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 259
$ws1.Range("F12").Value = 108
$ws1.Range("F22").Value = 1874
$ws1.Range("F23").Value = 4006
$ws1.Range("F28").Value = 2091
$ws1.Range("F35").Value = 460

# Sheet "全部类型" (same events, shifted by one row vs "展览")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 259
$ws4.Range("F12").Value = 108
$ws4.Range("F23").Value = 1874
$ws4.Range("F24").Value = 4006
$ws4.Range("F29").Value = 2091
$ws4.Range("F36").Value = 460
